# Natmi following Dr Hou advice
#
# This recomputes the NATMI ligand-receptor scoring table after the number of
# ligand/receptor expressing cells (columns E and K) changed from 1 to 3 for
# every sending/target cluster pair on the sheet. All of the downstream
# columns (total expression, derived specificities, edge weights, ...) are
# re-derived from that change using the same relationships NATMI itself uses:
#
#   H (Ligand total expr)   = G (Ligand avg expr)   * E (Ligand-expressing cells)
#   N (Receptor total expr) = M (Receptor avg expr)  * K (Receptor-expressing cells)
#   I = G / SUM(G over all sending clusters for the row's target cluster block)
#   J = H / SUM(H over all sending clusters for the row's target cluster block)
#   O = M / SUM(M over all target clusters)
#   P = N / SUM(N over all target clusters)
#   Q (Edge avg weight)   = G * M
#   R (Edge total weight) = H * N
#   S = Q / SUM(Q over every row)
#   T = R / SUM(R over every row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 17
$firstDataRow = 2

# New ligand/receptor-expressing cell counts (previously 1, now 3 everywhere).
$newExpressingCells = 3

# Per-sending-cluster ligand average expression values (column G), keyed by
# the row at which each sending cluster block starts.
$ligandAvgBySendingBlock = @{
    2  = 2.512729          # ECs   (rows 2-5)
    6  = 92.89399466666667 # FAPs  (rows 6-9)
    10 = 1.610639333333333 # M2    (rows 10-13)
    14 = 17.68108066666667 # sCs   (rows 14-17)
}

# Per-target-cluster receptor average expression values (column M), keyed by
# the offset of the row within each 4-row sending-cluster block (0-based).
$receptorAvgByTargetOffset = @{
    0 = 27.53580066666666  # target cluster ECs  (D=20)
    1 = 27.50472933333333  # target cluster FAPs (D=21)
    2 = 45.12975566666668  # target cluster M2   (D=22)
    3 = 33.831228          # target cluster sCs  (D=23)
}

# Update E (Ligand-expressing cells) and K (Receptor-expressing cells), and
# set the new ligand/receptor average expression values (G, M) which are the
# inputs that everything else below is derived from.
for ($blockStart = $firstDataRow; $blockStart -le $lastRow; $blockStart += 4) {
    $ligandAvg = $ligandAvgBySendingBlock[$blockStart]
    for ($offset = 0; $offset -le 3; $offset++) {
        $r = $blockStart + $offset
        if ($r -gt $lastRow) { continue }

        $ws.Range("E$r").Value = $newExpressingCells
        $ws.Range("K$r").Value = $newExpressingCells

        $ws.Range("G$r").Value = $ligandAvg
        $ws.Range("M$r").Value = $receptorAvgByTargetOffset[$offset]
    }
}

# Totals used for the specificity (derived-specificity) normalizations below.
$sumLigandAvg = 0
$sumLigandTotal = 0
foreach ($blockStart in $ligandAvgBySendingBlock.Keys) {
    $sumLigandAvg += $ligandAvgBySendingBlock[$blockStart]
    $sumLigandTotal += ($ligandAvgBySendingBlock[$blockStart] * $newExpressingCells)
}

$sumReceptorAvg = 0
$sumReceptorTotal = 0
foreach ($offset in $receptorAvgByTargetOffset.Keys) {
    $sumReceptorAvg += $receptorAvgByTargetOffset[$offset]
    $sumReceptorTotal += ($receptorAvgByTargetOffset[$offset] * $newExpressingCells)
}

# First pass: compute H, N, Q, R for every row and accumulate the grand
# totals needed for the edge-level specificity columns S and T.
$edgeAvgWeight = @{}
$edgeTotalWeight = @{}
$sumEdgeAvgWeight = 0
$sumEdgeTotalWeight = 0

for ($blockStart = $firstDataRow; $blockStart -le $lastRow; $blockStart += 4) {
    $ligandAvg = $ligandAvgBySendingBlock[$blockStart]
    $ligandTotal = $ligandAvg * $newExpressingCells

    for ($offset = 0; $offset -le 3; $offset++) {
        $r = $blockStart + $offset
        if ($r -gt $lastRow) { continue }

        $receptorAvg = $receptorAvgByTargetOffset[$offset]
        $receptorTotal = $receptorAvg * $newExpressingCells

        $q = $ligandAvg * $receptorAvg
        $rr = $ligandTotal * $receptorTotal

        $edgeAvgWeight[$r] = $q
        $edgeTotalWeight[$r] = $rr
        $sumEdgeAvgWeight += $q
        $sumEdgeTotalWeight += $rr
    }
}

# Second pass: write H, I, J, N, O, P, Q, R, S, T using the totals above.
for ($blockStart = $firstDataRow; $blockStart -le $lastRow; $blockStart += 4) {
    $ligandAvg = $ligandAvgBySendingBlock[$blockStart]
    $ligandTotal = $ligandAvg * $newExpressingCells

    for ($offset = 0; $offset -le 3; $offset++) {
        $r = $blockStart + $offset
        if ($r -gt $lastRow) { continue }

        $receptorAvg = $receptorAvgByTargetOffset[$offset]
        $receptorTotal = $receptorAvg * $newExpressingCells

        $ws.Range("H$r").Value = $ligandTotal
        $ws.Range("I$r").Value = ($ligandAvg / $sumLigandAvg)
        $ws.Range("J$r").Value = ($ligandTotal / $sumLigandTotal)

        $ws.Range("N$r").Value = $receptorTotal
        $ws.Range("O$r").Value = ($receptorAvg / $sumReceptorAvg)
        $ws.Range("P$r").Value = ($receptorTotal / $sumReceptorTotal)

        $ws.Range("Q$r").Value = $edgeAvgWeight[$r]
        $ws.Range("R$r").Value = $edgeTotalWeight[$r]
        $ws.Range("S$r").Value = ($edgeAvgWeight[$r] / $sumEdgeAvgWeight)
        $ws.Range("T$r").Value = ($edgeTotalWeight[$r] / $sumEdgeTotalWeight)
    }
}
